$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Second paragraph (currently empty) becomes a paragraph with a left tab
#    stop at 990 twips and a single tab-character run.
# ---------------------------------------------------------------------------
$pTab = $d.Paragraphs(2)
$rTab = $pTab.Range
$xmlTab = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="990"/></w:tabs></w:pPr><w:r><w:tab/></w:r></w:p><w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rTab.InsertXML($xmlTab)

# ---------------------------------------------------------------------------
# 2) Collapse the "Cada nivel tiene..." paragraph (previously split across
#    three runs around a gramStart/gramEnd proofErr pair) into one run with
#    the full sentence and no proofErr markers.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Cada nivel tiene una vista y una din", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pSentence = $d.Paragraphs(11)
$rSentence = $pSentence.Range
$xmlSentence = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Cada nivel tiene una vista y una dinámica distinta. Uno de los niveles estará controlado por tiempo.</w:t></w:r></w:p><w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rSentence.InsertXML($xmlSentence)

# ---------------------------------------------------------------------------
# 3) Re-anchor the automatic "_GoBack" bookmark so it spans the whole
#    document: start right before the very first character, end right after
#    the very last character (adding a bookmark with the reserved name
#    replaces/moves the existing one).
# ---------------------------------------------------------------------------
$docEnd = $d.Content.End
$fullRange = $d.Range(0, $docEnd - 1)
$d.Bookmarks.Add("_GoBack", $fullRange)

Write-Output "edit complete"
